$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the table by two rows, mirroring the formatting used by the row above.
$ws.Range("A36:F37").Copy()
$ws.Range("A38:F39").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 38: 四方坪站
$ws.Range("A38").Value = 45980
$ws.Range("B38").Value = "四方坪站"
$ws.Range("C38").Value = 9131.59
$ws.Range("D38").Value = 7940.23
$ws.Range("E38").Value = 3039.16
$ws.Range("F38").Value = 383

# Row 39: 高岭站
$ws.Range("A39").Value = 45980
$ws.Range("B39").Value = "高岭站"
$ws.Range("C39").Value = 4902.33
$ws.Range("D39").Value = 4410.03
$ws.Range("E39").Value = 1277.17
$ws.Range("F39").Value = 181

# Reflect the new selection state shown in the workbook after the edit.
$ws.Range("I40").Select()
